$wb = $excel.ActiveWorkbook

# --- "users" sheet: rename camelCase column headers to lowercase -----------
$wsUsers = $wb.Worksheets.Item("users")
$wsUsers.Range("C1").Value = "firstname"
$wsUsers.Range("D1").Value = "lastname"
$wsUsers.Range("E1").Value = "picname"
$wsUsers.Range("F1").Value = "acctype"

# --- "enrollments" sheet: rename courseCode -> coursecode, (re)bold header -
$wsEnrollments = $wb.Worksheets.Item("enrollments")
$wsEnrollments.Range("A1").Value = "coursecode"

# Bold the header row / key columns (mirrors the manual "ugly" re-bolding
# from the original edit; harmless no-op where it was already bold).
$wsEnrollments.Columns("A:B").Font.Bold = $true
$wsEnrollments.Range("A1:B1").Font.Bold = $true

# --- restore/update the per-sheet selections seen in the saved workbook ----
$wsUsers.Activate()
$wsUsers.Range("E8").Select()

$wsCourses = $wb.Worksheets.Item("courses")
$wsCourses.Activate()
$wsCourses.Range("E13").Select()

$wsEnrollments.Activate()
$wsEnrollments.Range("F7").Select()
